# New Submission Synced: 2026-02-08 20:21:25
# Sheet "JSS 3E" gets a new form-response row (row 10) appended, and the
# previous last row's "Admission No" (C9) gets normalized from a text value
# to a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3E")

# C9 was stored as text "35" -> becomes a genuine number 35
$ws.Range("C9").Value = 35

# Append the new submission as row 10
$ws.Range("A10").Value = "2026-02-08 20:21:25"
$ws.Range("B10").Value = "Shatu Musa Hassan "

# Admission No for this submission stays textual ("39"), so force the cell
# to Text before writing the value, then strip the formatting we added so
# the cell keeps the default (unstyled) look of its neighbours.
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "39"
$ws.Range("C10").ClearFormats()

$ws.Range("D10").Value = 10
